$wb = $excel.ActiveWorkbook

# Period_1
$ws1 = $wb.Worksheets.Item("Period_1")
$ws1.Range("C6").Value = 380
$ws1.Range("D6").Value = 379
$ws1.Range("C7").Value = 15
$ws1.Range("D7").Value = 25

# Period_2
$ws2 = $wb.Worksheets.Item("Period_2")
$ws2.Range("C6").Value = 380
$ws2.Range("D6").Value = 379
$ws2.Range("C7").Value = 15
$ws2.Range("D7").Value = 25

# Period_3
$ws3 = $wb.Worksheets.Item("Period_3")
$ws3.Range("C6").Value = 380
$ws3.Range("D6").Value = 379
$ws3.Range("C7").Value = 15
$ws3.Range("D7").Value = 25

# Update selections to match final state
$ws1.Activate()
$ws1.Range("C6:D7").Select()

$ws2.Activate()
$ws2.Range("C6:D7").Select()

$ws3.Activate()
$ws3.Range("E19").Select()
